$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "26.343.72"
$ws.Cells.Item(2, 5).Value = "  +0.67%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.622.45"
$ws.Cells.Item(3, 5).Value = "  +1.26%  "

$ws.Cells.Item(4, 5).Value = "  -0.06%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "212.19"
$ws.Cells.Item(5, 5).Value = "  +0.10%  "

$ws.Cells.Item(6, 5).Value = "  -0.08%  "

$ws.Cells.Item(8, 5).Value = "  +0.09%  "

$ws.Cells.Item(9, 5).Value = "  +0.41%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "18.81"
$ws.Cells.Item(10, 5).Value = "  +3.72%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0814"
$ws.Cells.Item(11, 5).Value = "  +0.33%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.847.90"
$ws.Cells.Item(12, 5).Value = "  +1.26%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.621.23"
$ws.Cells.Item(13, 5).Value = "  +1.13%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.04"
$ws.Cells.Item(14, 5).Value = "  +0.47%  "

$ws.Cells.Item(15, 5).Value = "  +0.63%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "26.362.59"
$ws.Cells.Item(16, 5).Value = "  +0.70%  "

$ws.Cells.Item(17, 5).Value = "  +2.50%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.0₃0726"
$ws.Cells.Item(18, 5).Value = "  +0.01%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "202.20"
$ws.Cells.Item(20, 5).Value = "  -0.87%  "

$ws.Cells.Item(21, 5).Value = "  -0.01%  "

$ws.Cells.Item(22, 5).Value = "  +0.39%  "

$ws.Cells.Item(23, 5).Value = "  +0.42%  "

$ws.Cells.Item(24, 5).Value = "  -3.44%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "144.61"
$ws.Cells.Item(25, 5).Value = "  +0.04%  "

$ws.Cells.Item(26, 5).Value = "  -0.09%  "

$ws.Cells.Item(27, 5).Value = "  -2.16%  "

$ws.Cells.Item(28, 2).Value = "Cosmos"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "6.60"
$ws.Cells.Item(28, 5).Value = "  +1.22%  "

$ws.Cells.Item(29, 2).Value = "EthereumClassic"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.19"
$ws.Cells.Item(29, 5).Value = "  -0.02%  "

$ws.Cells.Item(30, 5).Value = "  +5.45%  "

$ws.Cells.Item(31, 5).Value = "  +0.37%  "

$ws.Cells.Item(32, 5).Value = "  +1.59%  "

$ws.Cells.Item(33, 5).Value = "  +0.30%  "

$ws.Cells.Item(34, 5).Value = "  +0.49%  "

$ws.Cells.Item(35, 5).Value = "  +2.22%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.160.76"
$ws.Cells.Item(36, 5).Value = "  +1.65%  "

$ws.Cells.Item(37, 5).Value = "  +0.26%  "

$ws.Cells.Item(38, 5).Value = "  +2.06%  "

$ws.Cells.Item(39, 5).Value = "  -0.07%  "

$ws.Cells.Item(40, 5).Value = "  +0.01%  "

$ws.Cells.Item(41, 5).Value = "  -0.04%  "

$ws.Cells.Item(42, 5).Value = "  +3.95%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.783"
$ws.Cells.Item(43, 5).Value = "  -0.15%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.759.51"
$ws.Cells.Item(44, 5).Value = "  +1.22%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "92.24"
$ws.Cells.Item(45, 5).Value = "  +0.03%  "

$ws.Cells.Item(46, 5).Value = "  +9.73%  "

$ws.Cells.Item(47, 5).Value = "  +1.18%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "53.94"
$ws.Cells.Item(48, 5).Value = "  -0.35%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0507"
$ws.Cells.Item(49, 5).Value = "  +0.16%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.410"
$ws.Cells.Item(50, 5).Value = "  +0.94%  "

$ws.Cells.Item(51, 5).Value = "  -0.23%  "
